# Insert a new "Freelance" income row above the existing "Salary" row,
# and update the Salary row's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes the existing row 2 "Salary" down to row 3)
$ws.Rows.Item(2).Insert()

# New row 2: Freelance income
$ws.Cells.Item(2, 1).Value = "Freelance"
$ws.Cells.Item(2, 2).Value = 30000
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "16/5/2025"

# Row 3 (previously row 2 "Salary"): update the date
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = "1/4/2025"
